$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Define the named range "data" used by the MATCH/MODE/INDEX formulas below
$wb.Names.Add('data', '=Sheet1!$B$2:$B$9')

# Column C: for every subject in B2:B9, find the position of its first
# occurrence within the "data" range (this powers the "most frequent text"
# calculation below).
$ws.Range("C2").Formula = "=MATCH(data,data,0)"
$ws.Range("C3").Formula = "=MATCH(data,data,0)"
$ws.Range("C4").Formula = "=MATCH(data,data,0)"
$ws.Range("C5").Formula = "=MATCH(data,data,0)"
$ws.Range("C6").Formula = "=MATCH(data,data,0)"
$ws.Range("C7").Formula = "=MATCH(data,data,0)"
$ws.Range("C8").Formula = "=MATCH(data,data,0)"
$ws.Range("C9").Formula = "=MATCH(data,data,0)"

# F5: the most frequently occurring piece of text in the list, using the
# classic MATCH/MODE/INDEX "most frequent text" trick.
$ws.Range("F5").FormulaArray = "=INDEX(data,MODE(MATCH(data,data,0)))"

# Match the selection left behind in the saved file
[void]$ws.Range("E11").Select()
